$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '318.12'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '3.86%'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '39.83'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '2.37%'

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.141'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '0.99%'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.08213'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '2.00%'

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '2.061'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '6.23%'

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '8.314'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '3.91%'

$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.9342'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '0.23%'

$ws.Range('B9').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C9').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.1377'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '-4.85%'

$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1983'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '2.81%'

$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.09078'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '-0.48%'

$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.03484'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-0.71%'

$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.09813'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '0.33%'

$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.001413'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '1.38%'

$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.006335'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '5.44%'

$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.683'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '-2.75%'

$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.299'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '2.52%'

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.176'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '-7.50%'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.3470'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '1.49%'

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.1292'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-0.78%'

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.905'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '2.38%'

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.2451'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '1.50%'

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04327'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '-1.36%'

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.001229'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-0.78%'

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.004747'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '11.01%'

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0004001'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '-10.05%'

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02226'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '9.40%'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.05224'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '3.40%'

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007630'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '2.64%'

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.009675'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '-5.15%'

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1383'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '2.59%'

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.009200'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '0.91%'

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00006568'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '6.19%'

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.002778'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '-9.14%'
